# [FIX] corregir y mejorar test
# Fill in the missing "IVA ventas" (G2) value to match "IVA compras" (F2),
# and move the active selection to G3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2 was blank but already formatted as a percentage (same style as F2);
# give it the same 0.21 (21%) value.
$ws.Range("G2").Value = 0.21

# Update the active cell / selection.
$ws.Range("G3").Select() | Out-Null
